$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 values
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "Binary S 14"
$ws.Range("D6").Value = "Median of Array"
$ws.Range("E6").Value = "https://www.scaler.com/academy/mentee-dashboard/class/30364/homework/problems/198/?navref=cl_pb_nv_tb"
$ws.Range("F6").Value = "dsa/4_median_of_array.java at main · ankurnecessary/dsa · GitHub"

# Hyperlinks
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.scaler.com/academy/mentee-dashboard/class/30364/homework/problems/198/?navref=cl_pb_nv_tb")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/4_median_of_array.java", "", "dsa/4_median_of_array.java at main · ankurnecessary/dsa · GitHub")

# Row height for new row
$ws.Rows.Item(6).RowHeight = 72

# Style: B6, C6, D6 -> default style (index 2)
# E6 -> vertical center + wrap text hyperlink style
$ws.Range("E6").WrapText = $true
$ws.Range("E6").VerticalAlignment = -4108  # xlCenter

# F6 -> hyperlink style with wrap text only (like F4)
$ws.Range("F6").WrapText = $true

# Selection
$ws.Range("C7").Select()
